$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.88321066666667
$ws.Range("H2").Value = 107.649632
$ws.Range("I2").Value = 0.08317795499144418
$ws.Range("J2").Value = 0.08448843719082051
$ws.Range("M2").Value = 166.3936563333333
$ws.Range("N2").Value = 499.180969
$ws.Range("O2").Value = 0.6959913618211631
$ws.Range("P2").Value = 0.7009944564025758
$ws.Range("Q2").Value = 5970.738623805934
$ws.Range("R2").Value = 53736.64761425342
$ws.Range("S2").Value = 0.05789113816799464
$ws.Range("T2").Value = 0.0592259261008824
$ws.Range("G3").Value = 35.88321066666667
$ws.Range("H3").Value = 107.649632
$ws.Range("I3").Value = 0.08317795499144418
$ws.Range("J3").Value = 0.08448843719082051
$ws.Range("O3").Value = 0.2039972194837954
$ws.Range("P3").Value = 0.2054636419703505
$ws.Range("Q3").Value = 1750.041946402614
$ws.Range("R3").Value = 15750.37751762352
$ws.Range("S3").Value = 0.01696807154060289
$ws.Range("T3").Value = 0.01735930200960919
$ws.Range("G4").Value = 35.88321066666667
$ws.Range("H4").Value = 107.649632
$ws.Range("I4").Value = 0.08317795499144418
$ws.Range("J4").Value = 0.08448843719082051
$ws.Range("M4").Value = 7.402863
$ws.Range("N4").Value = 22.208589
$ws.Range("O4").Value = 0.03096469429353687
$ws.Range("P4").Value = 0.03118728224898178
$ws.Range("Q4").Value = 265.638492565472
$ws.Range("R4").Value = 2390.746433089248
$ws.Range("S4").Value = 0.002575579948271638
$ws.Range("T4").Value = 0.002634964737445489
$ws.Range("G5").Value = 35.88321066666667
$ws.Range("H5").Value = 107.649632
$ws.Range("I5").Value = 0.08317795499144418
$ws.Range("J5").Value = 0.08448843719082051
$ws.Range("M5").Value = 5.118919500000001
$ws.Range("N5").Value = 10.237839
$ws.Range("O5").Value = 0.02141141574965316
$ws.Range("P5").Value = 0.0143768870013594
$ws.Range("Q5").Value = 183.683266804208
$ws.Range("R5").Value = 1102.099600825248
$ws.Range("S5").Value = 0.001780957775527749
$ws.Range("T5").Value = 0.001214680714413878
$ws.Range("G6").Value = 35.88321066666667
$ws.Range("H6").Value = 107.649632
$ws.Range("I6").Value = 0.08317795499144418
$ws.Range("J6").Value = 0.08448843719082051
$ws.Range("M6").Value = 11.38837866666667
$ws.Range("N6").Value = 34.165136
$ws.Range("O6").Value = 0.04763530865185137
$ws.Range("P6").Value = 0.04797773237673265
$ws.Range("Q6").Value = 408.6515908477725
$ws.Range("R6").Value = 3677.864317629953
$ws.Range("S6").Value = 0.003962207559047244
$ws.Range("T6").Value = 0.004053563628469572
$ws.Range("I7").Value = 0.03522729558434242
$ws.Range("J7").Value = 0.03578230735158529
$ws.Range("M7").Value = 166.3936563333333
$ws.Range("N7").Value = 499.180969
$ws.Range("O7").Value = 0.6959913618211631
$ws.Range("P7").Value = 0.7009944564025758
$ws.Range("Q7").Value = 2528.710574566262
$ws.Range("R7").Value = 22758.39517109637
$ws.Range("S7").Value = 0.02451789342702312
$ws.Range("T7").Value = 0.02508319909075443
$ws.Range("I8").Value = 0.03522729558434242
$ws.Range("J8").Value = 0.03578230735158529
$ws.Range("O8").Value = 0.2039972194837954
$ws.Range("P8").Value = 0.2054636419703505
$ws.Range("S8").Value = 0.007186270349139638
$ws.Range("T8").Value = 0.007351963186559162
$ws.Range("I9").Value = 0.03522729558434242
$ws.Range("J9").Value = 0.03578230735158529
$ws.Range("M9").Value = 7.402863
$ws.Range("N9").Value = 22.208589
$ws.Range("O9").Value = 0.03096469429353687
$ws.Range("P9").Value = 0.03118728224898178
$ws.Range("Q9").Value = 112.502473728112
$ws.Range("R9").Value = 1012.522263553008
$ws.Range("S9").Value = 0.001090802438557224
$ws.Range("T9").Value = 0.001115952918893706
$ws.Range("I10").Value = 0.03522729558434242
$ws.Range("J10").Value = 0.03578230735158529
$ws.Range("M10").Value = 5.118919500000001
$ws.Range("N10").Value = 10.237839
$ws.Range("O10").Value = 0.02141141574965316
$ws.Range("P10").Value = 0.0143768870013594
$ws.Range("Q10").Value = 77.793025018168
$ws.Range("R10").Value = 466.758150109008
$ws.Range("S10").Value = 0.0007542662714922763
$ws.Range("T10").Value = 0.0005144381894416536
$ws.Range("I11").Value = 0.03522729558434242
$ws.Range("J11").Value = 0.03578230735158529
$ws.Range("M11").Value = 11.38837866666667
$ws.Range("N11").Value = 34.165136
$ws.Range("O11").Value = 0.04763530865185137
$ws.Range("P11").Value = 0.04797773237673265
$ws.Range("Q11").Value = 173.0709823689102
$ws.Range("R11").Value = 1557.638841320192
$ws.Range("S11").Value = 0.001678063098130152
$ws.Range("T11").Value = 0.001716753965936353
$ws.Range("G12").Value = 177.70077
$ws.Range("H12").Value = 533.10231
$ws.Range("I12").Value = 0.4119137160358794
$ws.Range("J12").Value = 0.4184034835782469
$ws.Range("M12").Value = 166.3936563333333
$ws.Range("N12").Value = 499.180969
$ws.Range("O12").Value = 0.6959913618211631
$ws.Range("P12").Value = 0.7009944564025758
$ws.Range("Q12").Value = 29568.28085354871
$ws.Range("R12").Value = 266114.5276819384
$ws.Range("S12").Value = 0.2866883881766276
$ws.Range("T12").Value = 0.2932985225278772
$ws.Range("G13").Value = 177.70077
$ws.Range("H13").Value = 533.10231
$ws.Range("I13").Value = 0.4119137160358794
$ws.Range("J13").Value = 0.4184034835782469
$ws.Range("O13").Value = 0.2039972194837954
$ws.Range("P13").Value = 0.2054636419703505
$ws.Range("Q13").Value = 8666.554514781152
$ws.Range("R13").Value = 77998.99063303035
$ws.Range("S13").Value = 0.08402925273855708
$ws.Range("T13").Value = 0.08596670354906834
$ws.Range("G14").Value = 177.70077
$ws.Range("H14").Value = 533.10231
$ws.Range("I14").Value = 0.4119137160358794
$ws.Range("J14").Value = 0.4184034835782469
$ws.Range("M14").Value = 7.402863
$ws.Range("N14").Value = 22.208589
$ws.Range("O14").Value = 0.03096469429353687
$ws.Range("P14").Value = 0.03118728224898178
$ws.Range("Q14").Value = 1315.49445530451
$ws.Range("R14").Value = 11839.45009774059
$ws.Range("S14").Value = 0.01275478229236576
$ws.Range("T14").Value = 0.013048867536312
$ws.Range("G15").Value = 177.70077
$ws.Range("H15").Value = 533.10231
$ws.Range("I15").Value = 0.4119137160358794
$ws.Range("J15").Value = 0.4184034835782469
$ws.Range("M15").Value = 5.118919500000001
$ws.Range("N15").Value = 10.237839
$ws.Range("O15").Value = 0.02141141574965316
$ws.Range("P15").Value = 0.0143768870013594
$ws.Range("Q15").Value = 909.6359367180152
$ws.Range("R15").Value = 5457.815620308091
$ws.Range("S15").Value = 0.008819655827028785
$ws.Range("T15").Value = 0.00601533960437959
$ws.Range("G16").Value = 177.70077
$ws.Range("H16").Value = 533.10231
$ws.Range("I16").Value = 0.4119137160358794
$ws.Range("J16").Value = 0.4184034835782469
$ws.Range("M16").Value = 11.38837866666667
$ws.Range("N16").Value = 34.165136
$ws.Range("O16").Value = 0.04763530865185137
$ws.Range("P16").Value = 0.04797773237673265
$ws.Range("Q16").Value = 2023.72365811824
$ws.Range("R16").Value = 18213.51292306416
$ws.Range("S16").Value = 0.01962163700130017
$ws.Range("T16").Value = 0.02007405036060978
$ws.Range("G17").Value = 20.074196
$ws.Range("H17").Value = 40.148392
$ws.Range("I17").Value = 0.04653236263856699
$ws.Range("J17").Value = 0.0315103250497358
$ws.Range("M17").Value = 166.3936563333333
$ws.Range("N17").Value = 499.180969
$ws.Range("O17").Value = 0.6959913618211631
$ws.Range("P17").Value = 0.7009944564025758
$ws.Range("Q17").Value = 3340.218870391975
$ws.Range("R17").Value = 20041.31322235185
$ws.Range("S17").Value = 0.03238612244157245
$ws.Range("T17").Value = 0.02208856317930801
$ws.Range("G18").Value = 20.074196
$ws.Range("H18").Value = 40.148392
$ws.Range("I18").Value = 0.04653236263856699
$ws.Range("J18").Value = 0.0315103250497358
$ws.Range("O18").Value = 0.2039972194837954
$ws.Range("P18").Value = 0.2054636419703505
$ws.Range("Q18").Value = 979.0284756470201
$ws.Range("R18").Value = 5874.17085388212
$ws.Range("S18").Value = 0.009492472594279313
$ws.Range("T18").Value = 0.006474226144388283
$ws.Range("G19").Value = 20.074196
$ws.Range("H19").Value = 40.148392
$ws.Range("I19").Value = 0.04653236263856699
$ws.Range("J19").Value = 0.0315103250497358
$ws.Range("M19").Value = 7.402863
$ws.Range("N19").Value = 22.208589
$ws.Range("O19").Value = 0.03096469429353687
$ws.Range("P19").Value = 0.03118728224898178
$ws.Range("Q19").Value = 148.606522823148
$ws.Range("R19").Value = 891.639136938888
$ws.Range("S19").Value = 0.001440860383859223
$ws.Range("T19").Value = 0.0009827214010832712
$ws.Range("G20").Value = 20.074196
$ws.Range("H20").Value = 40.148392
$ws.Range("I20").Value = 0.04653236263856699
$ws.Range("J20").Value = 0.0315103250497358
$ws.Range("M20").Value = 5.118919500000001
$ws.Range("N20").Value = 10.237839
$ws.Range("O20").Value = 0.02141141574965316
$ws.Range("P20").Value = 0.0143768870013594
$ws.Range("Q20").Value = 102.758193351222
$ws.Range("R20").Value = 411.0327734048881
$ws.Range("S20").Value = 0.0009963237622679854
$ws.Range("T20").Value = 0.0004530203826161562
$ws.Range("G21").Value = 20.074196
$ws.Range("H21").Value = 40.148392
$ws.Range("I21").Value = 0.04653236263856699
$ws.Range("J21").Value = 0.0315103250497358
$ws.Range("M21").Value = 11.38837866666667
$ws.Range("N21").Value = 34.165136
$ws.Range("O21").Value = 0.04763530865185137
$ws.Range("P21").Value = 0.04797773237673265
$ws.Range("Q21").Value = 228.6125454768854
$ws.Range("R21").Value = 1371.675272861312
$ws.Range("S21").Value = 0.002216583456588015
$ws.Range("T21").Value = 0.001511793942340079
$ws.Range("G22").Value = 182.547562
$ws.Range("H22").Value = 547.642686
$ws.Range("I22").Value = 0.423148670749767
$ws.Range("J22").Value = 0.4298154468296114
$ws.Range("M22").Value = 166.3936563333333
$ws.Range("N22").Value = 499.180969
$ws.Range("O22").Value = 0.6959913618211631
$ws.Range("P22").Value = 0.7009944564025758
$ws.Range("Q22").Value = 30374.75629591586
$ws.Range("R22").Value = 273372.8066632428
$ws.Range("S22").Value = 0.2945078196079453
$ws.Range("T22").Value = 0.3012982455037537
$ws.Range("G23").Value = 182.547562
$ws.Range("H23").Value = 547.642686
$ws.Range("I23").Value = 0.423148670749767
$ws.Range("J23").Value = 0.4298154468296114
$ws.Range("O23").Value = 0.2039972194837954
$ws.Range("P23").Value = 0.2054636419703505
$ws.Range("Q23").Value = 8902.934959783192
$ws.Range("R23").Value = 80126.41463804872
$ws.Range("S23").Value = 0.0863211522612165
$ws.Range("T23").Value = 0.0883114470807255
$ws.Range("G24").Value = 182.547562
$ws.Range("H24").Value = 547.642686
$ws.Range("I24").Value = 0.423148670749767
$ws.Range("J24").Value = 0.4298154468296114
$ws.Range("M24").Value = 7.402863
$ws.Range("N24").Value = 22.208589
$ws.Range("O24").Value = 0.03096469429353687
$ws.Range("P24").Value = 0.03118728224898178
$ws.Range("Q24").Value = 1351.374592470006
$ws.Range("R24").Value = 12162.37133223005
$ws.Range("S24").Value = 0.01310266923048302
$ws.Range("T24").Value = 0.01340477565524731
$ws.Range("G25").Value = 182.547562
$ws.Range("H25").Value = 547.642686
$ws.Range("I25").Value = 0.423148670749767
$ws.Range("J25").Value = 0.4298154468296114
$ws.Range("M25").Value = 5.118919500000001
$ws.Range("N25").Value = 10.237839
$ws.Range("O25").Value = 0.02141141574965316
$ws.Range("P25").Value = 0.0143768870013594
$ws.Range("Q25").Value = 934.4462747992591
$ws.Range("R25").Value = 5606.677648795555
$ws.Range("S25").Value = 0.009060212113336359
$ws.Range("T25").Value = 0.006179408110508124
$ws.Range("G26").Value = 182.547562
$ws.Range("H26").Value = 547.642686
$ws.Range("I26").Value = 0.423148670749767
$ws.Range("J26").Value = 0.4298154468296114
$ws.Range("M26").Value = 11.38837866666667
$ws.Range("N26").Value = 34.165136
$ws.Range("O26").Value = 0.04763530865185137
$ws.Range("P26").Value = 0.04797773237673265
$ws.Range("Q26").Value = 2078.920760732811
$ws.Range("R26").Value = 18710.2868465953
$ws.Range("S26").Value = 0.02015681753678578
$ws.Range("T26").Value = 0.02062157047937686
